$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.484.37"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "1.919.34"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4828"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4076"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08181"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("D12").Value = "1.911.68"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.041"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.226"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06795"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.008"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "29.517.70"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.641"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.197"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "2.145.37"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.677"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.022"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09555"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.521"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.384"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02285"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06132"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.182"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.039"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1857"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.284"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.399"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07608"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5574"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.959"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.432"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.38%  "
